# "minor changes fetching dd values from excel"
# Add a new "Java" entry below the existing "JavaScript" entry on the
# AutoCompleteSampleSheet worksheet (column A drop-down data list).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AutoCompleteSampleSheet")
$ws.Activate()

$ws.Range("A3").Value = "Java"
$ws.Range("A4").Select()
